$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns contain text that often *looks* numeric
# (e.g. "10.72", "1.00"). Force the range to Text format first so
# assigning .Value keeps these as strings instead of Excel silently
# recasting them to real numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.990.45"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "1.648.67"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "213.75"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "0.525"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "23.65"
$ws.Range("E8").Value = "  +4.14%  "
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "0.0872"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "1.880.98"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").Value = "1.653.68"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "0.563"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").Value = "65.72"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "27.939.36"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "231.95"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "10.72"
$ws.Range("E22").Value = "  +6.70%  "
$ws.Range("D23").Value = "4.39"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").Value = "152.04"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "1.20"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Value = "1.454.36"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "0.890"
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").Value = "0.562"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "0.917"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.02"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "2.24"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +6.04%  "
$ws.Range("D48").Value = "1.790.20"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "89.06"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.101"
$ws.Range("E51").Value = "  +1.43%  "

# Restore the default (unstyled) look so only cell values differ from
# the source workbook, not formatting.
$dataRange.Style = "Normal"
